# Fri, Jul 17, 2020  5:05:30 PM
# Renumber the embedded Pearson / BTEC logo pictures so their names no
# longer collide: the two Pearson logos (in the "first page" and
# "default" footers) become "image2.png" (they were sharing the name
# "image1.png" with the header's BTEC logo), and the BTEC logo (in the
# "first page" header) becomes "image1.jpg" (it was sharing the name
# "image2.jpg" with the main document picture).

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.Name -eq "image2.jpg") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.Name -eq "image1.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
